# Applies regenerated s_vals data (sam hentges, 2022) reflecting
# save-game filtering: updates columns B-E and the derived sum G
# for rows 2-13 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ "B"=3.286832544864788; "C"=250555.8564151394; "D"=3.537761648806719; "E"=1133.036916526867; "G"=251695.7179258599 }
    3 = @{ "B"=0.2917716402565462; "C"=1.655778082260271; "D"=0.7527432677738641; "E"=0.4942365360607697; "G"=3.194529526351451 }
    4 = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=0.7527432677738641; "E"=0.4942365360607697; "G"=6.189590430959694 }
    5 = @{ "B"=0.6606524410359556; "C"=1.655778082260271; "D"=0.1494219747398047; "E"=0.4942365360607697; "G"=2.960089034096801 }
    6 = @{ "B"=0.0006408296065709695; "C"=0.04071648406533734; "D"=22.3905356188092; "E"=10.19245300693656; "G"=32.62434593941767 }
    7 = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=3.537761648806719; "E"=0.4942365360607697; "G"=8.974608811992548 }
    8 = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=0.1494219747398047; "E"=0.4942365360607697; "G"=5.586269137925634 }
    9 = @{ "B"=1.455362044514542; "C"=1.655778082260271; "D"=3.537761648806719; "E"=0.4942365360607697; "G"=7.143138311642302 }
    10 = @{ "B"=1.455362044514542; "C"=1.655778082260271; "D"=0.7527432677738641; "E"=0.4942365360607697; "G"=4.358119930609447 }
    11 = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=22.3905356188092; "E"=0.4942365360607697; "G"=27.82738278199502 }
    12 = @{ "B"=0.6606524410359556; "C"=0.306821227259698; "D"=0.7527432677738641; "E"=0.4942365360607697; "G"=2.214453472130288 }
    13 = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=0.1494219747398047; "E"=0.4942365360607697; "G"=5.586269137925634 }
}

foreach ($row in $newValues.Keys) {
    foreach ($col in $newValues[$row].Keys) {
        $ws.Range("$col$row").Value = $newValues[$row][$col]
    }
}
